# Update Apoe-Lrp1 NATMI sheet with newly re-computed TPM-based values.
# The raw per-cluster ligand/receptor expression values (columns G,H for the
# "Sending cluster" and M,N for the "Target cluster") were recomputed from
# new TPM data. Every other numeric column on the sheet is derived from
# these via simple formulas, so we recompute them here as well:
#   I = G / SUM(G over all sending clusters)
#   J = H / SUM(H over all sending clusters)
#   O = M / SUM(M over all target clusters)
#   P = N / SUM(N over all target clusters)
#   Q = G * M
#   R = H * N
#   S = I * O
#   T = J * P

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ligand average/total expression value" per Sending cluster.
$newG = @{
    "ECs"               = 37.17328633333333
    "FAPs"              = 45.15892033333333
    "Inflammatory-Mac"  = 2375.59786
    "MuSCs"             = 76.954105
    "Resolving-Mac"     = 4655.195393666666
}
$newH = @{
    "ECs"               = 111.519859
    "FAPs"              = 135.476761
    "Inflammatory-Mac"  = 7126.79358
    "MuSCs"             = 153.90821
    "Resolving-Mac"     = 13965.586181
}

# New "Receptor average/total expression value" per Target cluster.
$newM = @{
    "ECs"               = 2.906846333333333
    "FAPs"              = 185.8027443333333
    "Inflammatory-Mac"  = 137.0717086666666
    "MuSCs"             = 29.2127365
    "Resolving-Mac"     = 171.5584106666666
}
$newN = @{
    "ECs"               = 8.720538999999999
    "FAPs"              = 557.408233
    "Inflammatory-Mac"  = 411.2151259999999
    "MuSCs"             = 58.425473
    "Resolving-Mac"     = 514.6752319999999
}

$firstRow = 2
$lastRow = 26

# First pass: write the updated raw values for columns G, H, M, N.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value()
    $target  = $ws.Cells.Item($r, 4).Value()

    $ws.Cells.Item($r, 7).Value  = $newG[$sending]
    $ws.Cells.Item($r, 8).Value  = $newH[$sending]
    $ws.Cells.Item($r, 13).Value = $newM[$target]
    $ws.Cells.Item($r, 14).Value = $newN[$target]
}

# Compute the totals needed for the specificity columns.
$sumG = 0.0
$sumH = 0.0
foreach ($v in $newG.Values) { $sumG += $v }
foreach ($v in $newH.Values) { $sumH += $v }

$sumM = 0.0
$sumN = 0.0
foreach ($v in $newM.Values) { $sumM += $v }
foreach ($v in $newN.Values) { $sumN += $v }

# Second pass: recompute the derived specificity and edge-weight columns.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value()
    $target  = $ws.Cells.Item($r, 4).Value()

    $g = $newG[$sending]
    $h = $newH[$sending]
    $m = $newM[$target]
    $n = $newN[$target]

    $i = $g / $sumG
    $j = $h / $sumH
    $o = $m / $sumM
    $p = $n / $sumN

    $ws.Cells.Item($r, 9).Value  = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p

    $ws.Cells.Item($r, 17).Value = $g * $m
    $ws.Cells.Item($r, 18).Value = $h * $n
    $ws.Cells.Item($r, 19).Value = $i * $o
    $ws.Cells.Item($r, 20).Value = $j * $p
}
